$wb = $excel.ActiveWorkbook

# The affected test data lives on the "Add Panels" worksheet.
$ws = $wb.Worksheets.Item("Add Panels")

# Fix the mislabeled Alarm/Standby current columns (rows 8-10):
#   Column H is the "AlarmLoadingDetail" column  -> should read "Alarm Current(A)"
#   Column I is the "StandbyLoadingDetail" column -> should read "Standby Current(A)"
# Update column I first, then column H, to mirror the authoring order of the change.
for ($r = 8; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "Standby Current(A)"
}
for ($r = 8; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = "Alarm Current(A)"
}

# Update the selected range on the sheet to match the new active selection.
$ws.Range("H9:H10").Select()
